$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 6471.4287
$ws.Range("I34").Value = 7750
$ws.Range("J34").Value = 5960
$ws.Range("K34").Value = 7750
$ws.Range("L34").Value = 5960
$ws.Range("M34").Value = -7547
$ws.Range("N34").Value = -6366
# Row 36
$ws.Range("H36").Value = 6471.4287
$ws.Range("I36").Value = 7750
$ws.Range("J36").Value = 5960
$ws.Range("K36").Value = 7750
$ws.Range("L36").Value = 5960
$ws.Range("M36").Value = -7035
$ws.Range("N36").Value = -7390
# Row 47
$ws.Range("H47").Value = 8000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
# Row 54
$ws.Range("H54").Value = 9974.5
$ws.Range("I54").Value = 9898
$ws.Range("K54").Value = 9898
$ws.Range("M54").Value = -9412
# Row 63
$ws.Range("H63").Value = 35333.332
$ws.Range("J63").Value = 35333.332
$ws.Range("L63").Value = 35333.332
$ws.Range("N63").Value = -36581.332
# Row 64
$ws.Range("H64").Value = 4184.4614
$ws.Range("I64").Value = 3857
$ws.Range("J64").Value = 4566.5
$ws.Range("K64").Value = 3857
$ws.Range("L64").Value = 4566.5
$ws.Range("M64").Value = -3609
$ws.Range("N64").Value = -5062.5
# Row 66
$ws.Range("H66").Value = 35333.332
$ws.Range("J66").Value = 35333.332
$ws.Range("L66").Value = 105999.996
$ws.Range("N66").Value = -112239.996
# Row 67
$ws.Range("H67").Value = 4184.4614
$ws.Range("I67").Value = 3857
$ws.Range("J67").Value = 4566.5
$ws.Range("K67").Value = 3857
$ws.Range("L67").Value = 4566.5
$ws.Range("M67").Value = -2999
$ws.Range("N67").Value = -6282.5
# Row 137
$ws.Range("H137").Value = 1887.6285
$ws.Range("I137").Value = 1929.25
$ws.Range("J137").Value = 1832.1333
$ws.Range("K137").Value = 5787.75
$ws.Range("L137").Value = 5496.3999
$ws.Range("M137").Value = -3237.75
$ws.Range("N137").Value = -10596.3999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3020.1555
$ws.Range("I32").Value = 1771.2285
$ws.Range("J32").Value = 7391.4
$ws.Range("K32").Value = 1771.2285
$ws.Range("L32").Value = 7391.4
$ws.Range("M32").Value = -1484.2285
$ws.Range("N32").Value = -7965.4
# Row 45
$ws.Range("H45").Value = 2026.6666
$ws.Range("I45").Value = 1469.8
$ws.Range("K45").Value = 1469.8
$ws.Range("M45").Value = -1092.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 13666.75
$ws.Range("J81").Value = 13666.75
$ws.Range("L81").Value = 13666.75
$ws.Range("N81").Value = -15788.75
# Row 84
$ws.Range("H84").Value = 13666.75
$ws.Range("J84").Value = 13666.75
$ws.Range("L84").Value = 41000.25
$ws.Range("N84").Value = -51608.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8659.112999999999
$ws.Range("I31").Value = 9787
$ws.Range("J31").Value = 5275.4546
$ws.Range("K31").Value = 9787
$ws.Range("L31").Value = 5275.4546
$ws.Range("M31").Value = -9492
$ws.Range("N31").Value = -5865.4546
# Row 34
$ws.Range("H34").Value = 8659.112999999999
$ws.Range("I34").Value = 9787
$ws.Range("J34").Value = 5275.4546
$ws.Range("K34").Value = 9787
$ws.Range("L34").Value = 5275.4546
$ws.Range("M34").Value = -9585
$ws.Range("N34").Value = -5679.4546
# Row 132
$ws.Range("H132").Value = 23791.959
$ws.Range("I132").Value = 39413.617
$ws.Range("K132").Value = 118240.851
$ws.Range("M132").Value = -115710.851
# Row 134
$ws.Range("H134").Value = 873.3333
$ws.Range("I134").Value = 841.6667
$ws.Range("K134").Value = 2525.0001
$ws.Range("M134").Value = 9.999899999999798

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 725
$ws.Range("I25").Value = 450
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 1350
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -1181
$ws.Range("N25").Value = -3338
# Row 30
$ws.Range("H30").Value = 725
$ws.Range("I30").Value = 450
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 1350
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = -1248
$ws.Range("N30").Value = -3204
# Row 62
$ws.Range("H62").Value = 7256.5
$ws.Range("J62").Value = 7256.5
$ws.Range("L62").Value = 21769.5
$ws.Range("N62").Value = -23141.5
# Row 63
$ws.Range("H63").Value = 4902.8
$ws.Range("I63").Value = 3333.3333
$ws.Range("J63").Value = 7257
$ws.Range("K63").Value = 9999.999899999999
$ws.Range("L63").Value = 21771
$ws.Range("M63").Value = -9250.999899999999
$ws.Range("N63").Value = -23269
# Row 64
$ws.Range("H64").Value = 2551.6
$ws.Range("I64").Value = 1500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4230
# Row 65
$ws.Range("H65").Value = 7256.5
$ws.Range("J65").Value = 7256.5
$ws.Range("L65").Value = 65308.5
$ws.Range("N65").Value = -72172.5
# Row 66
$ws.Range("H66").Value = 4902.8
$ws.Range("I66").Value = 3333.3333
$ws.Range("J66").Value = 7257
$ws.Range("K66").Value = 29999.9997
$ws.Range("L66").Value = 65313
$ws.Range("M66").Value = -26255.9997
$ws.Range("N66").Value = -72801
# Row 67
$ws.Range("H67").Value = 2551.6
$ws.Range("I67").Value = 1500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3564
# Row 68
$ws.Range("H68").Value = 1446.7858
$ws.Range("J68").Value = 1481.1538
$ws.Range("L68").Value = 4443.4614
$ws.Range("N68").Value = -6065.4614
# Row 71
$ws.Range("H71").Value = 1446.7858
$ws.Range("J71").Value = 1481.1538
$ws.Range("L71").Value = 13330.3842
$ws.Range("N71").Value = -21442.3842
# Row 107
$ws.Range("H107").Value = 3307.889
$ws.Range("I107").Value = 4779.6523
$ws.Range("J107").Value = 704
$ws.Range("K107").Value = 14338.9569
$ws.Range("L107").Value = 2112
$ws.Range("M107").Value = -12418.9569
$ws.Range("N107").Value = -5952
# Row 131
$ws.Range("H131").Value = 812.9299999999999
$ws.Range("J131").Value = 825.28864
$ws.Range("L131").Value = 2475.86592
$ws.Range("N131").Value = -12555.86592
# Row 134
$ws.Range("H134").Value = 6170.091
$ws.Range("I134").Value = 6337.35
$ws.Range("K134").Value = 19012.05
$ws.Range("M134").Value = -13942.05

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 23141.36
$ws.Range("I132").Value = 3140.889
$ws.Range("K132").Value = 9422.667000000001
$ws.Range("M132").Value = -6892.667000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1019.2222
$ws.Range("I22").Value = 440
$ws.Range("J22").Value = 1184.7142
$ws.Range("K22").Value = 440
$ws.Range("L22").Value = 1184.7142
$ws.Range("M22").Value = -145
$ws.Range("N22").Value = -1774.7142
# Row 27
$ws.Range("H27").Value = 1019.2222
$ws.Range("I27").Value = 440
$ws.Range("J27").Value = 1184.7142
$ws.Range("K27").Value = 440
$ws.Range("L27").Value = 1184.7142
$ws.Range("M27").Value = -333
$ws.Range("N27").Value = -1398.7142
# Row 122
$ws.Range("H122").Value = 3400.9
$ws.Range("I122").Value = 2660.8
$ws.Range("J122").Value = 4141
$ws.Range("K122").Value = 7982.400000000001
$ws.Range("L122").Value = 12423
$ws.Range("M122").Value = -5532.400000000001
$ws.Range("N122").Value = -17323
